# Add the "orario di accesso" placeholder into the "Tempi di accesso ai
# locali aziendali" line of the stage agreement template: the row of
# ellipsis placeholder characters is replaced with a single pair of
# ellipses wrapping the new {P_ORARIO_ACCESSO} merge field.

$d = $word.ActiveDocument

$ellipsis = [char]0x2026

$oldText = "Tempi di accesso ai locali aziendali: " + $ellipsis + $ellipsis + $ellipsis + $ellipsis + $ellipsis + $ellipsis + $ellipsis + $ellipsis + $ellipsis + $ellipsis + $ellipsis
$newText = "Tempi di accesso ai locali aziendali: " + $ellipsis + "{P_ORARIO_ACCESSO}" + $ellipsis

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
